$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make the header row bold (row 1, columns A:K)
$ws.Range("A1:K1").Font.Bold = $true

# Change the "target" column (G) value from "deuteron" to "d" for every data row
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 7).Value = "d"
}

# Update the active cell selection to match the saved state
$ws.Range("G17").Select() | Out-Null
